$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '19.935.64'
$ws.Range("E2").Value = '  -5.27%  '

$ws.Range("D3").Value = '1.412.76'
$ws.Range("E3").Value = '  -6.07%  '

$ws.Range("E4").Value = '  -0.61%  '

$ws.Range("D5").Value = '''1.001'
$ws.Range("E5").Value = '  -0.49%  '

$ws.Range("D6").Value = '''275.97'
$ws.Range("E6").Value = '  -2.77%  '

$ws.Range("D7").Value = '''0.3661'
$ws.Range("E7").Value = '  -4.53%  '

$ws.Range("D8").Value = '''0.3090'
$ws.Range("E8").Value = '  -1.00%  '

$ws.Range("D9").Value = '''39.72'
$ws.Range("E9").Value = '  -6.95%  '

$ws.Range("D10").Value = '''1.028'
$ws.Range("E10").Value = '  -0.61%  '

$ws.Range("D11").Value = '''0.06514'
$ws.Range("E11").Value = '  -6.39%  '

$ws.Range("E12").Value = '  -0.52%  '

$ws.Range("E13").Value = '  -1.69%  '

$ws.Range("D14").Value = '''17.53'
$ws.Range("E14").Value = '  -1.11%  '

$ws.Range("D15").Value = '''6.179'
$ws.Range("E15").Value = '  -2.79%  '

$ws.Range("D16").Value = '1.411.82'
$ws.Range("E16").Value = '  -6.87%  '

$ws.Range("D17").Value = '''0.00001016'
$ws.Range("E17").Value = '  -4.67%  '

$ws.Range("D18").Value = '''0.05661'
$ws.Range("E18").Value = '  -13.67%  '

$ws.Range("E19").Value = '  -0.49%  '

$ws.Range("D20").Value = '''71.04'
$ws.Range("E20").Value = '  -13.03%  '

$ws.Range("D21").Value = '''5.606'
$ws.Range("E21").Value = '  -5.99%  '

$ws.Range("D22").Value = '''14.67'
$ws.Range("E22").Value = '  -2.74%  '

$ws.Range("D23").Value = '''10.87'
$ws.Range("E23").Value = '  +0.19%  '

$ws.Range("D24").Value = '''2.235'
$ws.Range("E24").Value = '  -4.66%  '

$ws.Range("D25").Value = '19.951.93'
$ws.Range("E25").Value = '  -5.19%  '

$ws.Range("D26").Value = '''2.256'
$ws.Range("E26").Value = '  -3.42%  '

$ws.Range("D27").Value = '''132.83'
$ws.Range("E27").Value = '  -10.13%  '

$ws.Range("D28").Value = '''17.22'
$ws.Range("E28").Value = '  -4.13%  '

$ws.Range("D29").Value = '1.572.33'
$ws.Range("E29").Value = '  -6.47%  '

$ws.Range("D30").Value = '''109.49'
$ws.Range("E30").Value = '  -4.19%  '

$ws.Range("D31").Value = '''3.876'
$ws.Range("E31").Value = '  -19.21%  '

$ws.Range("D32").Value = '''5.246'
$ws.Range("E32").Value = '  -10.54%  '

$ws.Range("D33").Value = '''0.8144'
$ws.Range("E33").Value = '  -14.03%  '

$ws.Range("D34").Value = '''0.07676'
$ws.Range("E34").Value = '  -3.06%  '

$ws.Range("D35").Value = '''1.484'
$ws.Range("E35").Value = '  +0.83%  '

$ws.Range("D36").Value = '''8.302'
$ws.Range("E36").Value = '  -0.86%  '

$ws.Range("D37").Value = '''4.908'
$ws.Range("E37").Value = '  -2.74%  '

$ws.Range("D38").Value = '''0.05764'
$ws.Range("E38").Value = '  -0.06%  '

$ws.Range("D39").Value = '''0.9979'
$ws.Range("E39").Value = '  -0.76%  '

$ws.Range("D40").Value = '''0.02046'
$ws.Range("E40").Value = '  -3.32%  '

$ws.Range("D41").Value = '''10.40'
$ws.Range("E41").Value = '  -7.51%  '

$ws.Range("D42").Value = '''0.1881'
$ws.Range("E42").Value = '  -4.31%  '

$ws.Range("D43").Value = '''1.095'
$ws.Range("E43").Value = '  -5.48%  '

$ws.Range("B44").Value = 'TheSandbox'
$ws.Range("C44").Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range("D44").Value = '''0.5296'
$ws.Range("E44").Value = '  -5.27%  '

$ws.Range("B45").Value = 'EnergySwap'
$ws.Range("C45").Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range("D45").Value = '''12.41'
$ws.Range("E45").Value = '  -4.83%  '

$ws.Range("D46").Value = '''3.534'
$ws.Range("E46").Value = '  -3.47%  '

$ws.Range("D47").Value = '''0.5169'
$ws.Range("E47").Value = '  -4.50%  '

$ws.Range("D48").Value = '''115.44'
$ws.Range("E48").Value = '  +1.96%  '

$ws.Range("D49").Value = '''1.766'
$ws.Range("E49").Value = '  -3.80%  '

$ws.Range("D50").Value = '''1.030'
$ws.Range("E50").Value = '  -9.22%  '

$ws.Range("D51").Value = '''1.001'
$ws.Range("E51").Value = '  -0.52%  '
